$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "276.90"
Set-TextValue 2 5 "0.77%"
Set-TextValue 2 7 "7"

Set-TextValue 3 4 "27.15"
Set-TextValue 3 5 "1.26%"
Set-TextValue 3 7 "7"

Set-TextValue 4 4 "4.860"
Set-TextValue 4 5 "-0.03%"
Set-TextValue 4 7 "7"

Set-TextValue 5 4 "0.06412"
Set-TextValue 5 5 "1.35%"
Set-TextValue 5 7 "7"

Set-TextValue 6 4 "6.943"
Set-TextValue 6 5 "1.08%"
Set-TextValue 6 7 "7"

Set-TextValue 7 4 "1.179"
Set-TextValue 7 5 "-6.16%"
Set-TextValue 7 7 "7"

Set-TextValue 8 4 "0.8763"
Set-TextValue 8 5 "0.89%"
Set-TextValue 8 7 "7"

Set-TextValue 9 4 "0.1538"
Set-TextValue 9 5 "-1.17%"
Set-TextValue 9 7 "7"

Set-TextValue 10 4 "0.05126"
Set-TextValue 10 5 "2.52%"
Set-TextValue 10 7 "7"

Set-TextValue 11 4 "0.07485"
Set-TextValue 11 5 "0.35%"
Set-TextValue 11 7 "7"

Set-TextValue 12 4 "0.02956"
Set-TextValue 12 5 "0.14%"
Set-TextValue 12 7 "7"

Set-TextValue 13 5 "-0.15%"
Set-TextValue 13 7 "7"

Set-TextValue 14 4 "0.001570"
Set-TextValue 14 5 "-0.03%"
Set-TextValue 14 7 "7"

Set-TextValue 15 4 "0.0006376"
Set-TextValue 15 5 "1.35%"
Set-TextValue 15 7 "7"

Set-TextValue 16 4 "0.006123"
Set-TextValue 16 5 "2.46%"
Set-TextValue 16 7 "7"

Set-TextValue 17 4 "3.479"
Set-TextValue 17 5 "0.95%"
Set-TextValue 17 7 "7"

Set-TextValue 18 4 "3.308"
Set-TextValue 18 5 "-0.45%"
Set-TextValue 18 7 "7"

Set-TextValue 19 5 "0.13%"
Set-TextValue 19 7 "7"

Set-TextValue 20 5 "0.27%"
Set-TextValue 20 7 "7"

Set-TextValue 21 4 "0.1321"
Set-TextValue 21 5 "-0.97%"
Set-TextValue 21 7 "7"

Set-TextValue 22 4 "3.906"
Set-TextValue 22 5 "-0.39%"
Set-TextValue 22 7 "7"

Set-TextValue 23 4 "0.04414"
Set-TextValue 23 5 "1.42%"
Set-TextValue 23 7 "7"

Set-TextValue 24 7 "7"

Set-TextValue 25 5 "0.15%"
Set-TextValue 25 7 "7"

Set-TextValue 26 4 "0.003864"
Set-TextValue 26 5 "-9.00%"
Set-TextValue 26 7 "7"

Set-TextValue 27 5 "8.33%"
Set-TextValue 27 7 "7"

Set-TextValue 28 7 "7"

Set-TextValue 29 7 "7"

Set-TextValue 30 7 "7"

Set-TextValue 31 7 "7"

Set-TextValue 32 7 "7"

Set-TextValue 33 7 "7"

Set-TextValue 34 7 "7"

Set-TextValue 35 7 "7"

Set-TextValue 36 7 "7"

Set-TextValue 37 7 "7"

Set-TextValue 38 7 "7"

Set-TextValue 39 7 "7"

Set-TextValue 40 4 "0.04170"
Set-TextValue 40 5 "3.10%"
Set-TextValue 40 7 "7"

Set-TextValue 41 4 "0.006822"
Set-TextValue 41 5 "2.02%"
Set-TextValue 41 7 "7"

Set-TextValue 42 7 "7"

Set-TextValue 43 5 "-7.22%"
Set-TextValue 43 7 "7"

Set-TextValue 44 4 "0.01133"
Set-TextValue 44 5 "5.88%"
Set-TextValue 44 7 "7"

Set-TextValue 45 4 "0.00005313"
Set-TextValue 45 5 "0.51%"
Set-TextValue 45 7 "7"

Set-TextValue 46 7 "7"

Set-TextValue 47 5 "-7.42%"
Set-TextValue 47 7 "7"

Set-TextValue 48 7 "7"

Set-TextValue 49 7 "7"

Set-TextValue 50 7 "7"

Set-TextValue 51 7 "7"
